$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4461060000000001
$ws.Range("H2").Value = 1.338318
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 74.99480263236401
$ws.Range("R2").Value = 674.9532236912761
$ws.Range("S2").Value = 0.2984182258032519
$ws.Range("T2").Value = 0.298418225803252

# Row 3
$ws.Range("G3").Value = 0.4461060000000001
$ws.Range("H3").Value = 1.338318
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 72.71805976831401
$ws.Range("R3").Value = 654.462537914826
$ws.Range("S3").Value = 0.2893586437755394
$ws.Range("T3").Value = 0.2893586437755394

# Row 4
$ws.Range("G4").Value = 0.4461060000000001
$ws.Range("H4").Value = 1.338318
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 74.05070969418
$ws.Range("R4").Value = 666.4563872476201
$ws.Range("S4").Value = 0.294661504941043
$ws.Range("T4").Value = 0.294661504941043

# Row 5
$ws.Range("G5").Value = 0.4461060000000001
$ws.Range("H5").Value = 1.338318
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 29.54414354650601
$ws.Range("R5").Value = 265.897291918554
$ws.Range("S5").Value = 0.1175616254801657
$ws.Range("T5").Value = 0.1175616254801657
